$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.383.38'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.691.37'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.93'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5471'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +4.44%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06467'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.04'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07678'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.541'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.670.70'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5840'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008410'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.31'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').Value = '26.441.37'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.950'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.010'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.99'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.74'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.267'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.010'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.59'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1322'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +6.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.898'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.76'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06345'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.404'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +3.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.329'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.604'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.596'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6170'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.270'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = '1.124.38'
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01635'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8797'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.74'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.55'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000109'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -5.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.017'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.233'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.215'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +3.17%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05275'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('E51').Value = '  +0.17%  '
